# "Fruta / hortaliza, semanal" update
#
# A new weekly price record (2022-11-10 / serial 44875) is inserted right
# after the existing row for 2021-08-04 (row 274), pushing every
# subsequent data row down by one. The sheet's used range therefore grows
# from A1:R378 to A1:R379, and the record that used to be the last row
# (378, serial 44335) ends up as the new last row (379).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 275..378 down to 276..379, leaving a blank row 275.
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(275, 1).Value  = 3
$ws.Cells.Item(275, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(275, 3).Value  = "Coquimbo"
$ws.Cells.Item(275, 4).Value  = 44875
$ws.Cells.Item(275, 5).Value  = 5
$ws.Cells.Item(275, 6).Value  = 100112001
$ws.Cells.Item(275, 7).Value  = "Berenjena"
$ws.Cells.Item(275, 8).Value  = "Sin especificar"
$ws.Cells.Item(275, 9).Value  = "Primera"
$ws.Cells.Item(275, 10).Value = 50
$ws.Cells.Item(275, 11).Value = 12000
$ws.Cells.Item(275, 12).Value = 12000
$ws.Cells.Item(275, 13).Value = 12000
$ws.Cells.Item(275, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(275, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(275, 16).Value = 200
$ws.Cells.Item(275, 17).Value = 60
$ws.Cells.Item(275, 18).Value = "Hortaliza"
